$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.02, 1.048799459398664, 1.055990864809146, 1.060034219066393, 1.067807532965665, 1.045065753005537, 1.053841293123766, 1.058729803690407, 1.062762098412826, 1.07051440501086, 1.021744766994359),
    @(1.02, 1.049793138431244, 1.056798899371844, 1.061027246809139, 1.06880178176444, 1.045331985444127, 1.054483778315363, 1.059351402941292, 1.063569020930546, 1.071324066551917, 1.021962223635487),
    @(1.02, 1.05043636763193, 1.057322009595971, 1.061670905327272, 1.069445982873579, 1.045503214059783, 1.054899141576571, 1.059753222623563, 1.064091638580956, 1.071848208122231, 1.022102717157797),
    @(1.02, 1.050706841096351, 1.057541985869531, 1.061941762472051, 1.069717008883882, 1.045574948832458, 1.055073671638741, 1.059922051960148, 1.064311462266015, 1.072068612994248, 1.0221617287485),
    @(1.02, 1.050752258279186, 1.057578924386615, 1.061987255978124, 1.069762527269826, 1.04558697876066, 1.055102970786841, 1.059950393518595, 1.064348378362405, 1.072105623187234, 1.022171634002812),
    @(1.02, 1.050439981476839, 1.057324948692294, 1.06167452350112, 1.069449603536568, 1.045504173565096, 1.054901474004949, 1.059755478906177, 1.064094575424091, 1.071851152964884, 1.022103505877815),
    @(1.02, 1.049135225299693, 1.056263889892, 1.060369588423417, 1.068143366243082, 1.04515594293908, 1.054058499732644, 1.058939957757574, 1.063034700762064, 1.070787984461794, 1.021818301915566),
    @(1.02, 1.046838034352251, 1.054396186993613, 1.058078619002919, 1.06584820577455, 1.044534355738211, 1.052570284177329, 1.057499899165817, 1.061170816550918, 1.068916393474127, 1.021314099228411),
    @(1.02, 1.045307917750067, 1.053152466305932, 1.056557074237943, 1.064322590866009, 1.044114634296035, 1.051576303505508, 1.056537880081762, 1.059930797606386, 1.067669957264608, 1.020976881152113),
    @(1.02, 1.044645684560556, 1.05261426879783, 1.055899608522208, 1.063663059781448, 1.043931630517126, 1.051145470581477, 1.056120853861561, 1.059394475369513, 1.067130552739082, 1.020830608714613),
    @(1.02, 1.04439974970331, 1.052414410525403, 1.055655603700566, 1.063418242242034, 1.043863465536447, 1.050985375574421, 1.055965882418085, 1.059195354437597, 1.066930241285956, 1.020776238608556),
    @(1.02, 1.044452501385662, 1.052457278424597, 1.055707934171374, 1.06347074910518, 1.043878095714103, 1.051019719405923, 1.055999127446525, 1.059238062331581, 1.066973206615493, 1.020787902896586),
    @(1.02, 1.044625354526586, 1.052597747379685, 1.055879434775058, 1.063642819784129, 1.043925999836935, 1.05113223839043, 1.056108045284591, 1.059378014069057, 1.067113993960526, 1.020826115238523),
    @(1.02, 1.044731861427295, 1.052684301888941, 1.055985129443003, 1.063748859657833, 1.043955490082875, 1.051201556522773, 1.056175143989962, 1.059464255303152, 1.067200743994249, 1.020849654094905),
    @(1.02, 1.045351874702773, 1.05318819197927, 1.056600737081694, 1.064366384395696, 1.04412675310005, 1.051604887396342, 1.056565546981173, 1.059966404550932, 1.067705762373426, 1.020986583426944),
    @(1.02, 1.045740877871506, 1.053504361053541, 1.056987259777513, 1.064754028502344, 1.04423384431543, 1.051857770720724, 1.056810312240234, 1.060281554466599, 1.068022630735095, 1.021072407555643),
    @(1.02, 1.045967807438818, 1.053688810039569, 1.057212844360041, 1.064980237954655, 1.044296187059696, 1.052005231496683, 1.056953034708893, 1.060465435286003, 1.068207484605747, 1.021122442716786),
    @(1.02, 1.046045189736672, 1.053751707932658, 1.057289785308436, 1.065057387010468, 1.044317423663958, 1.052055504706165, 1.057001691736757, 1.060528143912785, 1.068270520050322, 1.02113949923277),
    @(1.02, 1.045699138345245, 1.053470435693056, 1.056945775845435, 1.064712427277122, 1.044222367032374, 1.05183064305176, 1.056784055890917, 1.060247735733782, 1.067988630658903, 1.021063201981517),
    @(1.02, 1.044574452255988, 1.052556381324046, 1.055828926394237, 1.063592144787935, 1.043911898487835, 1.051099106129488, 1.05607597360508, 1.059336799196487, 1.067072534281858, 1.020814863703343),
    @(1.02, 1.043867595567695, 1.051981981785513, 1.055127918173697, 1.062888714539632, 1.043715599929519, 1.050638787060225, 1.055630373311029, 1.058764595841032, 1.066496822768956, 1.020658503790457),
    @(1.02, 1.044242287193838, 1.052286452853414, 1.0554994219051, 1.063261527224394, 1.04381976516446, 1.050882846062723, 1.055866632248497, 1.059067880359338, 1.066801992053005, 1.020741413881633),
    @(1.02, 1.045717998545897, 1.053485765001514, 1.056964520239542, 1.064731224761257, 1.044227553498742, 1.05184290100679, 1.056795920143936, 1.060263016780533, 1.068003993737411, 1.021067361660242),
    @(1.02, 1.047431678509667, 1.054878787853519, 1.05866987708993, 1.066440772191073, 1.044695992001679, 1.052955350064763, 1.057872540841519, 1.061652225820374, 1.069400020620536, 1.021444639874366)
)

$cols = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowData[$j]
    }
}
